$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A13").Value = "Elizaveta"
$ws.Range("B13").Value = "Lizavainer"
$ws.Range("A14").Value = "tester"
$ws.Range("B14").Value = "test"
